$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D51").Value = "[VUE] 컴포넌트 이름은 multi-word이어야 한다고 합니다"
$ws.Range("E51").Value = "https://bskyvision.com/entry/VUE-%EC%BB%B4%ED%8F%AC%EB%84%8C%ED%8A%B8-%EC%9D%B4%EB%A6%84%EC%9D%80-multi-word%EC%9D%B4%EC%96%B4%EC%95%BC-%ED%95%9C%EB%8B%A4%EA%B3%A0-%ED%95%A9%EB%8B%88%EB%8B%A4"

$ws.Range("D52").Value = "숨은 DS"
